# Renames ontology identifiers/labels to PascalCase conventions and updates the
# Parameters section (adds Azimuth/Elevation rows, shifts Datum* rows down).
# Regenerated from the Google Sheet source -> new .ttl export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 'ontolidar:Units'
$ws.Range("A23").Value = 'ontolidar:VelocityAzimuthDisplay'
$ws.Range("B23").Value = 'Velocity-azimuth display'
$ws.Range("A25").Value = 'ontolidar:UseCase'
$ws.Range("B25").Value = 'Use case'
$ws.Range("A26").Value = 'ontolidar:VirtualMetMast'
$ws.Range("B26").Value = 'Virtual met mast'
$ws.Range("A28").Value = 'ontolidar:ChassisModule'
$ws.Range("A29").Value = 'ontolidar:ControlModule'
$ws.Range("A30").Value = 'ontolidar:CommunicationsModule'
$ws.Range("A31").Value = 'ontolidar:SignalProcessingModule'
$ws.Range("A32").Value = 'ontolidar:StorageModule'
$ws.Range("A33").Value = 'ontolidar:SafetyInterlocks'
$ws.Range("A34").Value = 'ontolidar:OpticsModule'
$ws.Range("A36").Value = 'ontolidar:TelescopeAperture'
$ws.Range("A37").Value = 'ontolidar:TelescopeApertureDiameter'
$ws.Range("A38").Value = 'ontolidar:PhotonicsModule'
$ws.Range("A40").Value = 'ontolidar:Photodetector'
$ws.Range("B40").Value = 'Photodetector'
$ws.Range("A41").Value = 'ontolidar:PhotodetectorGain'
$ws.Range("B41").Value = 'Photodetector gain'
$ws.Range("A42").Value = 'ontolidar:PhotodetectorVoltageNoise'
$ws.Range("A43").Value = 'ontolidar:OpticalAmplifier'
$ws.Range("A44").Value = 'ontolidar:BeamSplitter'
$ws.Range("A45").Value = 'ontolidar:LaserSource'
$ws.Range("A46").Value = 'ontolidar:LaserDiode'
$ws.Range("A47").Value = 'ontolidar:LaserDiodeWavelength'
$ws.Range("A48").Value = 'ontolidar:PowerModule'
$ws.Range("A49").Value = 'ontolidar:UninterruptiblePowerSupply'
$ws.Range("A50").Value = 'ontolidar:UpsBattery'
$ws.Range("A51").Value = 'ontolidar:UpsBatteryVoltage'
$ws.Range("A52").Value = 'ontolidar:UpsBatteryCapacity'
$ws.Range("A53").Value = 'ontolidar:ScannerModule'
$ws.Range("A54").Value = 'ontolidar:ScannerAzimuthPositioningServo'
$ws.Range("A55").Value = 'ontolidar:AzimuthSlewRate'
$ws.Range("A56").Value = 'ontolidar:ScannerElevationPositioningServo'
$ws.Range("A57").Value = 'ontolidar:ElevationSlewRate'
$ws.Range("A58").Value = 'ontolidar:ScannerMirrors'
$ws.Range("A59").Value = 'ontolidar:Devices'
$ws.Range("B59").Value = 'Devices'
$ws.Range("E61").Value = 'Type of free-standing lidar device intended for use as a forward-looking lidar mounted on a wind turbine nacelle'
$ws.Range("A62").Value = 'ontolidar:Instances'
$ws.Range("B62").Value = 'Instances'
$ws.Range("E62").Value = 'Wind lidar instances'
$ws.Range("A63").Value = 'ontolidar:SerialNumber'
$ws.Range("B63").Value = 'Serial number'
$ws.Range("A64").Value = 'ontolidar:LidarType'
$ws.Range("B64").Value = 'Lidar type'
$ws.Range("G64").Value = 'ontolidar:Instances'
$ws.Range("A65").Value = 'ontolidar:MeasurementPrinciples'
$ws.Range("A66").Value = 'ontolidar:LineOfSight'
$ws.Range("B66").Value = 'Line of sight'
$ws.Range("A67").Value = 'ontolidar:TimeOfFlight'
$ws.Range("B67").Value = 'Time of flight'
$ws.Range("A68").Value = 'ontolidar:MeasurementVolume'
$ws.Range("A70").Value = 'ontolidar:ProbeVolume'
$ws.Range("A71").Value = 'ontolidar:RadialVelocity'
$ws.Range("B71").Value = 'Radial velocity'
$ws.Range("A72").Value = 'ontolidar:ScanningGeometry'
$ws.Range("F72").Value = 'Scan pattern'
$ws.Range("A73").Value = 'ontolidar:SingleLidar'
$ws.Range("G73").Value = 'ontolidar:ScanningGeometry'
$ws.Range("A74").Value = 'ontolidar:Staring'
$ws.Range("B74").Value = 'Staring'
$ws.Range("G74").Value = 'ontolidar:SingleLidar'
$ws.Range("A75").Value = 'ontolidar:StopStare'
$ws.Range("B75").Value = 'Stop-stare'
$ws.Range("G75").Value = 'ontolidar:SingleLidar'
$ws.Range("A76").Value = 'ontolidar:SweepStare'
$ws.Range("B76").Value = 'Sweep-stare'
$ws.Range("G76").Value = 'ontolidar:SingleLidar'
$ws.Range("A77").Value = 'ontolidar:SimpleScan'
$ws.Range("G77").Value = 'ontolidar:SingleLidar'
$ws.Range("A78").Value = 'ontolidar:VariableAzimuthScan'
$ws.Range("G78").Value = 'ontolidar:SimpleScan'
$ws.Range("A79").Value = 'ontolidar:CompleteCone'
$ws.Range("G79").Value = 'ontolidar:VariableAzimuthScan'
$ws.Range("A80").Value = 'ontolidar:ConeSector'
$ws.Range("G80").Value = 'ontolidar:VariableAzimuthScan'
$ws.Range("A81").Value = 'ontolidar:DopplerBeamSwinging'
$ws.Range("G81").Value = 'ontolidar:VariableAzimuthScan'
$ws.Range("A82").Value = 'ontolidar:VariableElevation'
$ws.Range("G82").Value = 'ontolidar:SimpleScan'
$ws.Range("A83").Value = 'ontolidar:VerticalSlice'
$ws.Range("G83").Value = 'ontolidar:VariableElevationScan'
$ws.Range("A84").Value = 'ontolidar:CompoundScan'
$ws.Range("G84").Value = 'ontolidar:SingleLidar'
$ws.Range("A85").Value = 'ontolidar:ArbitraryTrajectory'
$ws.Range("G85").Value = 'ontolidar:SingleLidar'
$ws.Range("A86").Value = 'ontolidar:SequentialScan'
$ws.Range("G86").Value = 'ontolidar:CompoundScan'
$ws.Range("A87").Value = 'ontolidar:MultiLidar'
$ws.Range("G87").Value = 'ontolidar:ScanningGeometry'
$ws.Range("A88").Value = 'ontolidar:ConvergentScan'
$ws.Range("G88").Value = 'ontolidar:MultiLidar'
$ws.Range("A89").Value = 'ontolidar:Parameters'
$ws.Range("B89").Value = 'Parameters'
$ws.Range("A90").Value = 'ontolidar:CarrierToNoiseRatio'
$ws.Range("G90").Value = 'ontolidar:Parameters'
$ws.Range("A91").Value = 'ontolidar:Azimuth'
$ws.Range("B91").Value = 'Azimuth'
$ws.Range("E91").Value = 'The angle between the line of sight of the lidar and a reference vector on the datum plane.'
$ws.Range("F91").Value = 'Azimuth angle'
$ws.Range("G91").Value = 'ontolidar:Parameters'
$ws.Range("H91").Value = 'degrees'
$ws.Range("M91").Value = 'The reference vector from which the azimuth is defined is often true north, but could be some other defined vector.'
$ws.Range("A92").Value = 'ontolidar:Elevation'
$ws.Range("B92").Value = 'Elevation'
$ws.Range("E92").Value = 'The angle between the line of sight of the lidar and the datum plane.'
$ws.Range("F92").Value = 'Elevation angle'
$ws.Range("G92").Value = 'ontolidar:Parameters'
$ws.Range("H92").Value = 'degrees'
$ws.Range("M92").Value = ''
$ws.Range("A93").Value = 'ontolidar:MeasurementHeight'
$ws.Range("B93").Value = 'Measurement height'
$ws.Range("E93").Value = 'The nominal height above the datum plane at which a windfield reconstruction process returns a wind speed. Often used for vertically-profiling wind lidars for comparison to point wind speed measurements from an anemometer.'
$ws.Range("F93").Value = ''
$ws.Range("G93").Value = 'ontolidar:Parameters'
$ws.Range("H93").Value = 'meters'
$ws.Range("I93").Value = ''
$ws.Range("A94").Value = 'ontolidar:DatumElevation'
$ws.Range("B94").Value = 'Datum elevation'
$ws.Range("E94").Value = 'The height of the datum plane above sea level'
$ws.Range("F94").Value = ''
$ws.Range("G94").Value = 'ontolidar:Parameters'
$ws.Range("H94").Value = 'meters'
$ws.Range("M94").Value = 'Sea level should be defined on a project basis and is out of scope of this schema.'
$ws.Range("A95").Value = 'ontolidar:DatumPlane'
$ws.Range("B95").Value = 'Datum plane'
$ws.Range("E95").Value = 'The horizontal plane from which the measurement height is defined, e.g., lidar window, local ground, a platform top surface, or other reference.'
$ws.Range("F95").Value = 'Reference height'
$ws.Range("G95").Value = 'ontolidar:Parameters'
$ws.Range("I95").Value = 'ontolidar:datumfeature'
$ws.Range("A96").Value = 'ontolidar:DatumFeature'
$ws.Range("B96").Value = 'Datum feature'
$ws.Range("E96").Value = 'A distinguishing feature used to recognise or define the datum plane from which the measurement height is defined'
$ws.Range("F96").Value = 'Reference marker'
$ws.Range("G96").Value = 'ontolidar:Parameters'
